$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("sNo"), shifting sourceCity..class right by one
$ws.Columns.Item(1).Insert()

# Copy the formatting of the (now shifted) neighboring header/data cells onto the
# new column so the new cells reuse the existing header/data styles
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# Header for the new first column
$ws.Range("A1").Value = "sNo"

# Row 2 gets the current dataset number
$ws.Range("A2").Value = 1

# Reset the selection to the top-left cell
$ws.Range("A1").Select() | Out-Null
